# Add the new hydrological cycle "15-16" (row 19) below the existing data
# and extend the yearly averages (row 23) so they include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New data row (cycle 15-16) --------------------------------------
$ws.Range("A19").Value2 = "15-16"

# Match the look of the other "Ciclo" labels in column A (A6:A18).
$ws.Range("A6").Copy()
$ws.Range("A19").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B19").Value2 = 4.6754859999999994
$ws.Range("C19").Value2 = 13.309269
$ws.Range("D19").Value2 = 6.7845400000000007
$ws.Range("E19").Value2 = 3.19693
$ws.Range("F19").Value2 = 21.835500000000003
$ws.Range("G19").Value2 = 10.226281999999999
$ws.Range("H19").Value2 = 42.415078999999999
$ws.Range("I19").Value2 = 117.12362299999999
$ws.Range("J19").Value2 = 170.40609000000001
$ws.Range("K19").Value2 = 202.70077000000001
$ws.Range("L19").Value2 = 80.746818000000005
$ws.Range("M19").Value2 = 31.789432000000001

# Formatting picked up from the pasted-in source data (Arial 10, 2 decimals).
$monthly = $ws.Range("B19:M19")
$monthly.Font.Name = "Arial"
$monthly.Font.Size = 10
$monthly.NumberFormat = "0.00"

$ws.Range("N19").Formula = "=SUM(B19:M19)"

# ---- Extend the "Promedio" row so it covers the new cycle too --------
$ws.Range("B23").Formula = "=AVERAGE(B6:B19)"
$ws.Range("C23:M23").Formula = "=AVERAGE(C6:C19)"

# ---- Match the saved selection in the edited file ---------------------
$ws.Range("N23").Select()
